# Move nodegoat files so spreadsheet folder (#82)
# Update the "Directory" column (D) for all data rows (3-20) from
# "data/nodegoat/" to "data/spreadsheets/".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($r = 3; $r -le 20; $r++) {
    $ws.Cells.Item($r, 4).Value = "data/spreadsheets/"
}

# Reflect the new selection captured in the saved workbook.
$ws.Range("G19").Select()
